# Update cosinor analysis results for sawtooth_10 after re-running CircadiPy simulations
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("E2").Value = 25.8200000000006
$ws.Range("G2").Value = [double]"1.984045309155746e-07"
$ws.Range("H2").Value = [double]"4.446137295974263e-06"
$ws.Range("K2").Value = 5.767194287030375
$ws.Range("L2").Value = "[3.369417638768116, 8.164970935292633]"
$ws.Range("M2").Value = [double]"3.301162605495733e-06"
$ws.Range("N2").Value = [double]"6.602325210991467e-06"
$ws.Range("O2").Value = -1.345947603339772
$ws.Range("P2").Value = "[-1.8365266363327732, -0.8553685703467702]"
$ws.Range("Q2").Value = [double]"1.289861939479664e-07"
$ws.Range("R2").Value = [double]"1.289861939479664e-07"
$ws.Range("S2").Value = 10.30005448423578
$ws.Range("T2").Value = "[8.890838667507182, 11.70927030096437]"
$ws.Range("W2").Value = 5.531011011011142
$ws.Range("X2").Value = 3.515035035035116
$ws.Range("Y2").Value = 7.546986986987168

# --- Row 3 ---
$ws.Range("E3").Value = 23.99000000000031
$ws.Range("G3").Value = [double]"3.056597588368959e-06"
$ws.Range("H3").Value = [double]"1.799523062150863e-05"
$ws.Range("K3").Value = 5.813530009289712
$ws.Range("L3").Value = "[2.7527402090029014, 8.874319809576523]"
$ws.Range("M3").Value = 0.0002236703294962172
$ws.Range("N3").Value = 0.0002236703294962172
$ws.Range("O3").Value = 2.119553001521041
$ws.Range("P3").Value = "[1.5912371198362703, 2.6478688832058115]"
$ws.Range("Q3").Value = [double]"6.23945339839338e-14"
$ws.Range("R3").Value = [double]"1.247890679678676e-13"
$ws.Range("S3").Value = 10.64006568984861
$ws.Range("T3").Value = "[9.070342437503236, 12.209788942193992]"
$ws.Range("W3").Value = 15.89727727727749
$ws.Range("X3").Value = 13.88010010010028
$ws.Range("Y3").Value = 17.91445445445469
